$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the bold/bordered header style introduced previously: revert header row to default formatting
$ws.Range("A1:L1").ClearFormats()

# Rewrite rows 2-9 with the updated collection data (new row order + refreshed temporal extents)
# Row 2: cems-glofas-forecast
$ws.Range("A2").Value = 'Collection'
$ws.Range("B2").Value = '1.1.0'
$ws.Range("C2").Value = 'cems-glofas-forecast'
$ws.Range("D2").Value = 'River discharge and related forecasted data by the Global Flood Awareness System'
$ws.Range("E2").Value = 'This dataset provides an ensemble of forecast time series of gridded hydrological data. The data set is a product of the Global Flood Awareness System (GloFAS) and offers a consistent representation of key hydrological variables across the global domain including:

River discharge
Soil wetness index (root zone)
Snow water equivalent
Runoff water equivalent (surface plus subsurface)

Also provided are two ancillary files for interpretation, one containing upstream area data and the other containing elevation data (see the table of related variables and the associated link in the documentation).
This dataset was produced by forcing the open-source LISFLOOD hydrological model with input from the European Centre for Medium-range Weather Forecasts (ECMWF) ensemble forecast combined with the ECMWF extended-range ensemble forecast up to 30 days.
Companion datasets, also available through the EWDS, are historical simulations which can be used to derive the hydrological climatology and for verification; reforecasts for research, local skill assessment and post-processing; and seasonal forecasts and reforecasts for users looking for longer leadtime forecasts. For users looking specifically for European hydrological data, we refer to the European Flood Awareness System (EFAS) forecasts and historical simulations.
All the GloFAS and EFAS datasets are part of the operational flood forecasting within the Copernicus Emergency Management Service (CEMS), which is managed, technically implemented and developed by the European Commission’s Joint Research Centre.'
$ws.Range("F2").Value = '[{''rel'': ''self'', ''type'': ''application/json'', ''href'': ''https://ewds.climate.copernicus.eu/api/catalogue/v1/collections/cems-glofas-forecast''}, {''rel'': ''parent'', ''type'': ''application/json'', ''href'': ''https://ewds.climate.copernicus.eu/api/catalogue/v1/''}, {''rel'': ''root'', ''type'': ''application/json'', ''href'': ''https://ewds.climate.copernicus.eu/api/catalogue/v1/''}]'
$ws.Range("G2").Value = '[''Temporal coverage: Past'', ''Spatial coverage: Global'', ''Provider: Copernicus CEMS'', ''Variable domain: Land (hydrology)'', ''Product type: Forecasts'']'
$ws.Range("H2").Value = 'other'
$ws.Range("I2").Value = '[{''name'': ''ECMWF''}]'
$ws.Range("J2").Value = '{''spatial'': {''bbox'': [[0, -70, 360, 70]]}, ''temporal'': {''interval'': [[''2019-11-05T00:00:00Z'', ''2026-02-03T00:00:00Z'']]}}'
$ws.Range("K2").Value = '{}'
$ws.Range("L2").Value = '{''thumbnail'': {''href'': ''https://object-store.os-api.cci2.ecmwf.int:443/cci2-prod-catalogue/resources/cems-glofas-forecast/overview_36fc7b601512e3619bc5ba70ae0488b911d9d74e203400f9a321f5745768f6a5.png'', ''type'': ''image/jpg'', ''roles'': [''thumbnail''], ''thumborHref'': ''https://ewds.climate.copernicus.eu/thumbnails/5Ei37SA8oZby9C0Xvw2VoF6dWag=/150x0/filters:format(webp)/object-store.os-api.cci2.ecmwf.int/cci2-prod-catalogue/resources/cems-glofas-forecast/overview_36fc7b601512e3619bc5ba70ae0488b911d9d74e203400f9a321f5745768f6a5.png''}}'

# Row 3: efas-historical
$ws.Range("A3").Value = 'Collection'
$ws.Range("B3").Value = '1.1.0'
$ws.Range("C3").Value = 'efas-historical'
$ws.Range("D3").Value = 'River discharge and related historical data from the European Flood Awareness System'
$ws.Range("E3").Value = 'This dataset provides gridded modelled sub-daily and daily hydrological time series forced with meteorological observations. The data set is a consistent representation of the most important hydrological variables across the European Flood Awareness System (EFAS) domain. The temporal resolution is up to 30 years modelled time series of:

River discharge
Volumetric soil moisture
Snow water equivalent
Soil wetness index (root zone)
Runoff water equivalent (surface plus subsurface)

Also provided are auxiliary (time invariant) data to aid interpretation of river discharge and soil moisture data. These auxiliary data are the upstream area, elevation, soil depth, wilting capacity and field capacity. The latter three are provided at three soil levels, one for each of the three soil layers represented in LISFLOOD.
This dataset was produced by forcing the open-source LISFLOOD hydrological model with gridded observational data of precipitation and temperature at a 1x1 arcminute resolution (~1.5 km at EFAS latitudes) across the EFAS domain. Previous versions of the data have a 5x5km resolution. For the latest version data is available from 1992-01-01 up until near-real time, with a delay of 6 days. The real-time data is only available to EFAS partners.
Companion datasets, also available through the EWDS, are forecasts for users who are looking medium-range forecasts, reforecasts for research, local skill assessment and post-processing, and seasonal forecasts and reforecasts for users looking for long-term forecasts. For users looking for global hydrological data, we refer to the Global Flood Awareness System (GloFAS) forecasts and historical simulations.
All these datasets are part of the operational flood forecasting within the Copernicus Emergency Management Service (CEMS), which is managed, technically implemented and developed by the European Commission’s Joint Research Centre.'
$ws.Range("F3").Value = '[{''rel'': ''self'', ''type'': ''application/json'', ''href'': ''https://ewds.climate.copernicus.eu/api/catalogue/v1/collections/efas-historical''}, {''rel'': ''parent'', ''type'': ''application/json'', ''href'': ''https://ewds.climate.copernicus.eu/api/catalogue/v1/''}, {''rel'': ''root'', ''type'': ''application/json'', ''href'': ''https://ewds.climate.copernicus.eu/api/catalogue/v1/''}]'
$ws.Range("G3").Value = '[''Product type: Reanalysis'', ''Temporal coverage: Past'', ''Spatial coverage: Europe'', ''Provider: Copernicus CEMS'', ''Variable domain: Land (hydrology)'']'
$ws.Range("H3").Value = 'other'
$ws.Range("I3").Value = '[{''name'': ''ECMWF''}]'
$ws.Range("J3").Value = '{''spatial'': {''bbox'': [[0, -70, 360, 70]]}, ''temporal'': {''interval'': [[''1991-01-01T00:00:00Z'', ''2026-01-31T00:00:00Z'']]}}'
$ws.Range("K3").Value = '{}'
$ws.Range("L3").Value = '{''thumbnail'': {''href'': ''https://object-store.os-api.cci2.ecmwf.int:443/cci2-prod-catalogue/resources/efas-historical/overview_ddd9074d456be00a54d03c320485bdbb1d1871507eccaa1039404a9c2c62fe31.png'', ''type'': ''image/jpg'', ''roles'': [''thumbnail''], ''thumborHref'': ''https://ewds.climate.copernicus.eu/thumbnails/lHpkSZaxXko-wA_Tt4GqwZhv-7U=/150x0/filters:format(webp)/object-store.os-api.cci2.ecmwf.int/cci2-prod-catalogue/resources/efas-historical/overview_ddd9074d456be00a54d03c320485bdbb1d1871507eccaa1039404a9c2c62fe31.png''}}'

# Row 4: cems-fire-historical-v1
$ws.Range("A4").Value = 'Collection'
$ws.Range("B4").Value = '1.1.0'
$ws.Range("C4").Value = 'cems-fire-historical-v1'
$ws.Range("D4").Value = 'Fire danger indices historical data from the Copernicus Emergency Management Service'
$ws.Range("E4").Value = 'This data set provides complete historical reconstruction of meteorological conditions favourable to the start, spread and sustainability of fires. The fire danger metrics  provided are part of a vast dataset produced by the Copernicus Emergency Management Service for the
European Forest Fire Information System (EFFIS). The European Forest Fire Information System incorporates the fire danger indices for three different models developed in Canada, United States and Australia. In this dataset the fire danger indices are calculated using weather forecast from historical simulations provided by ECMWF ERA5 reanalysis.
ERA5 by combining  model data and  a vast set of quality controlled observations provides a  globally complete and consistent data-set and is regarded as a good proxy for observed atmospheric conditions.
The selected data records in this data set are regularly extended with time as ERA5 forcing data become available. 
This dataset is produced by ECMWF in its role of the computational centre for fire danger forecast of the CEMS,  on behalf of the Joint Research Centre which is the managing entity of the service.'
$ws.Range("F4").Value = '[{''rel'': ''self'', ''type'': ''application/json'', ''href'': ''https://ewds.climate.copernicus.eu/api/catalogue/v1/collections/cems-fire-historical-v1''}, {''rel'': ''parent'', ''type'': ''application/json'', ''href'': ''https://ewds.climate.copernicus.eu/api/catalogue/v1/''}, {''rel'': ''root'', ''type'': ''application/json'', ''href'': ''https://ewds.climate.copernicus.eu/api/catalogue/v1/''}]'
$ws.Range("G4").Value = '[''Product type: Reanalysis'', ''Temporal coverage: Past'', ''Spatial coverage: Global'', ''Variable domain: Land (biosphere)'', ''Provider: Copernicus CEMS'']'
$ws.Range("H4").Value = 'CC-BY-4.0'
$ws.Range("I4").Value = '[{''name'': ''Copernicus Emergency Management Service''}]'
$ws.Range("J4").Value = '{''spatial'': {''bbox'': [[0, -90, 360, 90]]}, ''temporal'': {''interval'': [[''1940-01-03T00:00:00Z'', ''2026-01-31T00:00:00Z'']]}}'
$ws.Range("K4").Value = '{}'
$ws.Range("L4").Value = '{''thumbnail'': {''href'': ''https://object-store.os-api.cci2.ecmwf.int:443/cci2-prod-catalogue/resources/cems-fire-historical-v1/overview_7d52aa9bfd619151f13dba1d6a0c625edcb157696f3fbf0db3c24f0993f45e77.png'', ''type'': ''image/jpg'', ''roles'': [''thumbnail''], ''thumborHref'': ''https://ewds.climate.copernicus.eu/thumbnails/ObVmZLmujWun1v8ODWaaRNfQEAo=/150x0/filters:format(webp)/object-store.os-api.cci2.ecmwf.int/cci2-prod-catalogue/resources/cems-fire-historical-v1/overview_7d52aa9bfd619151f13dba1d6a0c625edcb157696f3fbf0db3c24f0993f45e77.png''}}'

# Row 5: cems-glofas-historical
$ws.Range("A5").Value = 'Collection'
$ws.Range("B5").Value = '1.1.0'
$ws.Range("C5").Value = 'cems-glofas-historical'
$ws.Range("D5").Value = 'River discharge and related historical data from the Global Flood Awareness System'
$ws.Range("E5").Value = 'This dataset provides gridded modelled daily hydrological time series forced with meteorological reanalysis data. The data set is a product of the Global Flood Awareness System (GloFAS) and offers a consistent representation of key hydrological variables across the global domain including:

River discharge
Soil wetness index (root zone)
Snow water equivalent
Runoff water equivalent (surface plus subsurface)

Also provided are two ancillary files for interpretation, one containing upstream area data and the other containing elevation data (see the table of related variables and the associated link in the documentation).
This dataset was produced by forcing the open-source LISFLOOD hydrological model with ERA5 meteorological reanalysis data, interpolated to the GloFAS resolution, produced at a 24-hourly timestep. Two variations of the ERA5 forcing data are used, resulting in two types of hydrological data: intermediate and consolidated. Intermediate hydrological data is produced using ERA5 Near Real Time (ERA5T) data and is updated daily, whilst consolidated hydrological data is produced using the consolidated ERA5 reanalysis and is updated monthly.
Companion datasets, also available through the EWDS, are forecasts for users who are looking for medium-range forecasts, reforecasts for research, local skill assessment and post-processing, and seasonal forecasts and reforecasts for users looking for long-term forecasts. For users specifically interested in European hydrological data, we refer to the European Flood Awareness System (EFAS) forecasts and historical simulations. 
All the GloFAS and EFAS datasets are part of the operational flood forecasting within the Copernicus Emergency Management Service (CEMS), which is managed, technically implemented and developed by the European Commission’s Joint Research Centre.'
$ws.Range("F5").Value = '[{''rel'': ''self'', ''type'': ''application/json'', ''href'': ''https://ewds.climate.copernicus.eu/api/catalogue/v1/collections/cems-glofas-historical''}, {''rel'': ''parent'', ''type'': ''application/json'', ''href'': ''https://ewds.climate.copernicus.eu/api/catalogue/v1/''}, {''rel'': ''root'', ''type'': ''application/json'', ''href'': ''https://ewds.climate.copernicus.eu/api/catalogue/v1/''}]'
$ws.Range("G5").Value = '[''Product type: Reanalysis'', ''Temporal coverage: Past'', ''Spatial coverage: Global'', ''Provider: Copernicus CEMS'', ''Variable domain: Land (hydrology)'']'
$ws.Range("H5").Value = 'other'
$ws.Range("I5").Value = '[{''name'': ''ECMWF''}]'
$ws.Range("J5").Value = '{''spatial'': {''bbox'': [[0, -70, 360, 70]]}, ''temporal'': {''interval'': [[''1979-01-01T00:00:00Z'', ''2026-01-31T00:00:00Z'']]}}'
$ws.Range("K5").Value = '{}'
$ws.Range("L5").Value = '{''thumbnail'': {''href'': ''https://object-store.os-api.cci2.ecmwf.int:443/cci2-prod-catalogue/resources/cems-glofas-historical/overview_779934fbef02d194554ae626a1c5570e24a2e2a2e4b2ae7e8c54b4ca0c7fe9d0.png'', ''type'': ''image/jpg'', ''roles'': [''thumbnail''], ''thumborHref'': ''https://ewds.climate.copernicus.eu/thumbnails/326qPafXpnwSyVeWE_zprZ7Bl9Y=/150x0/filters:format(webp)/object-store.os-api.cci2.ecmwf.int/cci2-prod-catalogue/resources/cems-glofas-historical/overview_779934fbef02d194554ae626a1c5570e24a2e2a2e4b2ae7e8c54b4ca0c7fe9d0.png''}}'

# Row 6: efas-forecast
$ws.Range("A6").Value = 'Collection'
$ws.Range("B6").Value = '1.1.0'
$ws.Range("C6").Value = 'efas-forecast'
$ws.Range("D6").Value = 'River discharge and related forecasted data by the European Flood Awareness System'
$ws.Range("E6").Value = 'This dataset provides gridded modelled hydrological time series forced with medium-range meteorological forecasts. The data is a consistent representation of the most important hydrological variables across the European Flood Awareness System (EFAS) domain. The temporal resolution is sub-daily high-resolution and ensemble forecasts of:

River discharge
Volumetric soil moisture
Snow water equivalent
Soil wetness index (root zone)
Runoff water equivalent (surface plus subsurface)

Also provided are auxiliary (time invariant) data to aid interpretation of river discharge and soil moisture data. These auxiliary data are the upstream area, elevation, soil depth, wilting capacity and field capacity. The latter three are provided at three soil levels, one for each of the three soil layers represented in LISFLOOD.
This data set was produced by forcing the open-source LISFLOOD hydrological model at a 1x1 arcminute resolution (~1.5 km at EFAS latitudes) with meteorological forecasts. For version 4.0 and older, the open-source LISFLOOD hydrological model was forced at a 5x5km resolution. The forecasts are initialised twice daily at 00 and 12 UTC with time steps of 6 or 24 hours (the latter for version 3.5 and older) and lead times between 5 and 15 days depending on the forcing numerical weather prediction model. The forcing meteorological data are high-resolution and ensemble forecasts from the European Centre of Medium-range Weather Forecasts (ECMWF) with 51 ensemble members, high-resolution forecasts from the Deutsches Wetter Dienst (DWD) and the ensemble forecasts from the COSMO Local Ensemble Prediction System (COSMO-LEPS) with 20 ensemble members. The hydrological forecasts are available from 2018-10-10 up until present with a 30-day delay. The real-time data is only available to EFAS partners.
Companion datasets, also available through the EWDS, are historical simulations which can be used to derive the hydrological climatology and for verification; reforecasts for research, local skill assessment and post-processing; and seasonal forecasts and reforecasts for users looking for longer leadtime forecasts. For users looking for global hydrological data, we refer to the Global Flood Awareness System (GloFAS) forecasts and historical simulations.
All these datasets are part of the operational flood forecasting within the Copernicus Emergency Management Service (CEMS), which is managed, technically implemented and developed by the European Commission’s Joint Research Centre.'
$ws.Range("F6").Value = '[{''rel'': ''self'', ''type'': ''application/json'', ''href'': ''https://ewds.climate.copernicus.eu/api/catalogue/v1/collections/efas-forecast''}, {''rel'': ''parent'', ''type'': ''application/json'', ''href'': ''https://ewds.climate.copernicus.eu/api/catalogue/v1/''}, {''rel'': ''root'', ''type'': ''application/json'', ''href'': ''https://ewds.climate.copernicus.eu/api/catalogue/v1/''}]'
$ws.Range("G6").Value = '[''Temporal coverage: Past'', ''Spatial coverage: Europe'', ''Provider: Copernicus CEMS'', ''Variable domain: Land (hydrology)'', ''Product type: Forecasts'']'
$ws.Range("H6").Value = 'other'
$ws.Range("I6").Value = '[{''name'': ''ECMWF''}]'
$ws.Range("J6").Value = '{''spatial'': {''bbox'': [[0, -70, 360, 70]]}, ''temporal'': {''interval'': [[''2018-10-11T00:00:00Z'', ''2025-12-29T00:00:00Z'']]}}'
$ws.Range("K6").Value = '{}'
$ws.Range("L6").Value = '{''thumbnail'': {''href'': ''https://object-store.os-api.cci2.ecmwf.int:443/cci2-prod-catalogue/resources/efas-forecast/overview_ddd9074d456be00a54d03c320485bdbb1d1871507eccaa1039404a9c2c62fe31.png'', ''type'': ''image/jpg'', ''roles'': [''thumbnail''], ''thumborHref'': ''https://ewds.climate.copernicus.eu/thumbnails/qaKs1FDYFXDR1n4PompbSKLekhY=/150x0/filters:format(webp)/object-store.os-api.cci2.ecmwf.int/cci2-prod-catalogue/resources/efas-forecast/overview_ddd9074d456be00a54d03c320485bdbb1d1871507eccaa1039404a9c2c62fe31.png''}}'

# Row 7: cems-glofas-seasonal
$ws.Range("A7").Value = 'Collection'
$ws.Range("B7").Value = '1.1.0'
$ws.Range("C7").Value = 'cems-glofas-seasonal'
$ws.Range("D7").Value = 'Seasonal forecasts of river discharge and related data by the Global Flood Awareness System'
$ws.Range("E7").Value = 'This dataset provides an ensemble of forecasted time series of gridded hydrological data forced with seasonal meteorological forecasts. The data set is a product of the Global Flood Awareness System (GloFAS) and offers a consistent representation of key hydrological variables across the global domain including:

River discharge
Soil wetness index (root zone)
Snow water equivalent
Runoff water equivalent (surface plus subsurface)

Also provided are two ancillary files for interpretation, one containing upstream area data and the other containing elevation data (see the table of related variables and the associated link in the documentation).
This dataset was produced by forcing the open-source LISFLOOD hydrological model at a 0.05° (~5 km at the equator) resolution with seasonal meteorological forecasts from the European Centre for Medium-range Weather Forecasts (ECMWF). For version 3.1 and older, the open-source LISFLOOD hydrological model was forced at a 0.1° (~11 km at the equator) resolution. The forecasts are initialised on the first of each month with a 24-hourly time step, and cover 123 days.
Companion datasets, also available through the Early Warning Data Store (EWDS), are the operational forecasts, historical simulations that can be used to derive the hydrological climatology, and medium-range and seasonal reforecasts. The latter dataset enables research, local skill assessment and post-processing of the seasonal forecasts. In addition, the seasonal reforecasts are also used to derive a specific range dependent climatology for the seasonal system. For users looking specifically for European hydrological data, we refer to the European Flood Awareness System (EFAS) forecasts and historical simulations.
All these datasets are part of the operational flood forecasting within the Copernicus Emergency Management Service (CEMS), which is managed, technically implemented and developed by the European Commission’s Joint Research Centre.'
$ws.Range("F7").Value = '[{''rel'': ''self'', ''type'': ''application/json'', ''href'': ''https://ewds.climate.copernicus.eu/api/catalogue/v1/collections/cems-glofas-seasonal''}, {''rel'': ''parent'', ''type'': ''application/json'', ''href'': ''https://ewds.climate.copernicus.eu/api/catalogue/v1/''}, {''rel'': ''root'', ''type'': ''application/json'', ''href'': ''https://ewds.climate.copernicus.eu/api/catalogue/v1/''}]'
$ws.Range("G7").Value = '[''Spatial coverage: Global'', ''Provider: Copernicus CEMS'', ''Product type: Seasonal forecasts'', ''Variable domain: Land (hydrology)'', ''Temporal coverage: Present'']'
$ws.Range("H7").Value = 'other'
$ws.Range("I7").Value = '[{''name'': ''ECMWF''}]'
$ws.Range("J7").Value = '{''spatial'': {''bbox'': [[0, -60, 360, 90]]}, ''temporal'': {''interval'': [[''2020-12-01T00:00:00Z'', ''2026-01-01T00:00:00Z'']]}}'
$ws.Range("K7").Value = '{}'
$ws.Range("L7").Value = '{''thumbnail'': {''href'': ''https://object-store.os-api.cci2.ecmwf.int:443/cci2-prod-catalogue/resources/cems-glofas-seasonal/overview_0885f764fd7ce14e5c511c5751d22b0610d3004d3de56f389562cace8f67e2bc.png'', ''type'': ''image/jpg'', ''roles'': [''thumbnail''], ''thumborHref'': ''https://ewds.climate.copernicus.eu/thumbnails/S_cANS_r0oitqGUjgW8BjGlgg1s=/150x0/filters:format(webp)/object-store.os-api.cci2.ecmwf.int/cci2-prod-catalogue/resources/cems-glofas-seasonal/overview_0885f764fd7ce14e5c511c5751d22b0610d3004d3de56f389562cace8f67e2bc.png''}}'

# Row 8: efas-seasonal
$ws.Range("A8").Value = 'Collection'
$ws.Range("B8").Value = '1.1.0'
$ws.Range("C8").Value = 'efas-seasonal'
$ws.Range("D8").Value = 'Seasonal forecasts of river discharge and related data by the European Flood Awareness System'
$ws.Range("E8").Value = 'This dataset provides gridded modelled daily hydrological time series forced with seasonal meteorological forecasts. The dataset is a consistent representation of the most important hydrological variables across the European Flood Awareness (EFAS) domain. The temporal resolution is daily forecasts initialised once a month consisting of:

River discharge
Volumetric soil moisture
Snow water equivalent
Soil wetness index (root zone)
Runoff water equivalent (surface plus subsurface)

Also provided are auxiliary (time invariant) data to aid interpretation of river discharge and soil moisture data. These auxiliary data are the upstream area, elevation, soil depth, wilting capacity and field capacity. The latter three are provided at three soil levels, one for each of the three soil layers represented in LISFLOOD.
This data set was produced by forcing the open-source LISFLOOD hydrological model at a 1x1 arcminute resolution (~1.5 km at EFAS latitudes) with seasonal meteorological ensemble forecasts. For version 4.0 and older, the open-source LISFLOOD hydrological model was forced at a 5x5km resolution. The forecasts are initialised on the first of each month with a lead time of 215 days at 24-hour time steps. The meteorological data are seasonal forecasts (SEAS5) from the European Centre of Medium-range Weather Forecasts (ECMWF) with 51 ensemble members. The forecasts are available from November 2020.
Companion datasets, also available through the Early Warning Data Store (EWDS), are seasonal reforecasts for research, local skill assessment and post-processing of the seasonal forecasts. There are also medium-range forecasts for users who want to look at shorter time ranges. These are accompanied by historical simulations which can be used to derive the hydrological climatology, and medium-range reforecasts. For users looking for global hydrological data, we refer to the Global Flood Awareness System (GloFAS) forecasts and historical simulations.
All these datasets are part of the operational flood forecasting within the Copernicus Emergency Management Service (CEMS), which is managed, technically implemented and developed by the European Commission’s Joint Research Centre.'
$ws.Range("F8").Value = '[{''rel'': ''self'', ''type'': ''application/json'', ''href'': ''https://ewds.climate.copernicus.eu/api/catalogue/v1/collections/efas-seasonal''}, {''rel'': ''parent'', ''type'': ''application/json'', ''href'': ''https://ewds.climate.copernicus.eu/api/catalogue/v1/''}, {''rel'': ''root'', ''type'': ''application/json'', ''href'': ''https://ewds.climate.copernicus.eu/api/catalogue/v1/''}]'
$ws.Range("G8").Value = '[''Spatial coverage: Europe'', ''Provider: Copernicus CEMS'', ''Product type: Seasonal forecasts'', ''Variable domain: Land (hydrology)'', ''Temporal coverage: Present'']'
$ws.Range("H8").Value = 'other'
$ws.Range("I8").Value = '[{''name'': ''ECMWF''}]'
$ws.Range("J8").Value = '{''spatial'': {''bbox'': [[0, -70, 360, 70]]}, ''temporal'': {''interval'': [[''2020-11-01T00:00:00Z'', ''2026-01-01T00:00:00Z'']]}}'
$ws.Range("K8").Value = '{}'
$ws.Range("L8").Value = '{''thumbnail'': {''href'': ''https://object-store.os-api.cci2.ecmwf.int:443/cci2-prod-catalogue/resources/efas-seasonal/overview_027e11ba4b15ebcb3217e766cb5c74cfd08665cf129d439a9ef52724a1961fa0.png'', ''type'': ''image/jpg'', ''roles'': [''thumbnail''], ''thumborHref'': ''https://ewds.climate.copernicus.eu/thumbnails/O1bh_HpiA2OYGzlcy0IGnZu_pAI=/150x0/filters:format(webp)/object-store.os-api.cci2.ecmwf.int/cci2-prod-catalogue/resources/efas-seasonal/overview_027e11ba4b15ebcb3217e766cb5c74cfd08665cf129d439a9ef52724a1961fa0.png''}}'

# Row 9: cems-fire-seasonal
$ws.Range("A9").Value = 'Collection'
$ws.Range("B9").Value = '1.1.0'
$ws.Range("C9").Value = 'cems-fire-seasonal'
$ws.Range("D9").Value = 'Seasonal forecast of fire danger indices from the Copernicus Emergency Management Service'
$ws.Range("E9").Value = 'This dataset offers modeled daily fire danger time series, driven by seasonal weather forecasts. It provides long-range predictions of meteorological conditions conducive to the initiation, spread, and persistence of fires. The fire danger metrics included in this dataset are part of an extensive dataset produced by the Copernicus Emergency Management Service (CEMS) for the European Forest Fire Information System (EFFIS) and the Global Wildfire Information System (GWIS). EFFIS and GWIS are used for monitoring and forecasting fire danger at both European and global scales. The dataset incorporates fire danger indices from the U.S. Forest Service National Fire-Danger Rating System (NFDRS), the Canadian Forest Service Fire Weather Index Rating System (FWI), and the Australian McArthur (Mark 5) rating systems.
This dataset was generated by driving the Global ECMWF Fire Forecast (GEFF) model with seasonal weather ensemble forecasts from the European Centre for Medium-Range Weather Forecasts (ECMWF) System 5 (SEAS5) prediction system.These forecasts initially consist of 25 ensemble members until December 2016, referred to as re-forecasts. After that period, they consist of seasonal forecasts with 51 members. It is important to note that the re-forecast dataset was initialized using ERA-Interim analysis data, while forecast simulations from 2016 onward are initialized using ECMWF operational analysis. Therefore, it is suggested that the period 1981-2016 be used as a reference period, while the period 2017-to present as a real time forecast. 
For both the re-forecast (1981-2016) and forecast periods (2017-present), the temporal resolution is daily forecasts at 12:00 local time, available once a month, with a prediction horizon of 216 days (equivalent to 7 months). The data records in this dataset will be extended over time as seasonal  forcing data becomes available. Once the SEAS5 operation ceases, the dataset will be updated with the next ECMWF seasonal system (SEAS6). It is essential to note that this is not a real-time service, as real-time forecasts are accessible through the EFFIS web services.
These seasonal forecasts can be used to assess the performance of the forecasting system or to develop tools for statistically correcting forecast errors. ECMWF produces this dataset as the computational center for fire danger forecasting within the Copernicus Emergency Management Service (CEMS) on behalf of the Joint Research Centre, which serves as the managing entity for this service.'
$ws.Range("F9").Value = '[{''rel'': ''self'', ''type'': ''application/json'', ''href'': ''https://ewds.climate.copernicus.eu/api/catalogue/v1/collections/cems-fire-seasonal''}, {''rel'': ''parent'', ''type'': ''application/json'', ''href'': ''https://ewds.climate.copernicus.eu/api/catalogue/v1/''}, {''rel'': ''root'', ''type'': ''application/json'', ''href'': ''https://ewds.climate.copernicus.eu/api/catalogue/v1/''}]'
$ws.Range("G9").Value = '[''Temporal coverage: Past'', ''Spatial coverage: Global'', ''Variable domain: Land (biosphere)'', ''Provider: Copernicus CEMS'', ''Product type: Seasonal forecasts'']'
$ws.Range("H9").Value = 'CC-BY-4.0'
$ws.Range("I9").Value = '[{''name'': ''ECMWF''}]'
$ws.Range("J9").Value = '{''spatial'': {''bbox'': [[-179.95, -89.95, 179.95, 89.95]]}, ''temporal'': {''interval'': [[''1981-02-01T00:00:00Z'', ''2025-12-01T00:00:00Z'']]}}'
$ws.Range("K9").Value = '{}'
$ws.Range("L9").Value = '{''thumbnail'': {''href'': ''https://object-store.os-api.cci2.ecmwf.int:443/cci2-prod-catalogue/resources/cems-fire-seasonal/overview_6a09913a8735343a5dc4801f5394a4448c819060f84fd021b7b4291e4ef231da.png'', ''type'': ''image/jpg'', ''roles'': [''thumbnail''], ''thumborHref'': ''https://ewds.climate.copernicus.eu/thumbnails/c8wgzziXyVvvElQIy6hl9nTQSrc=/150x0/filters:format(webp)/object-store.os-api.cci2.ecmwf.int/cci2-prod-catalogue/resources/cems-fire-seasonal/overview_6a09913a8735343a5dc4801f5394a4448c819060f84fd021b7b4291e4ef231da.png''}}'

